$d = $word.ActiveDocument

$replacements = @(
    @("2025-08-20 Wednesday", "2025-08-21 Thursday"),
    @("10÷3=3, 1", "44÷9=4, 8"),
    @("72÷2=36, 0", "21÷7=3, 0"),
    @("88÷6=14, 4", "82÷9=9, 1"),
    @("62÷5=12, 2", "29÷5=5, 4"),
    @("44÷8=5, 4", "94÷6=15, 4"),
    @("11÷8=1, 3", "25÷4=6, 1"),
    @("49÷8=6, 1", "10÷7=1, 3"),
    @("91÷8=11, 3", "73÷2=36, 1"),
    @("53÷4=13, 1", "54÷5=10, 4"),
    @("74÷8=9, 2", "84÷8=10, 4"),
    @("87÷8=10, 7", "35÷9=3, 8"),
    @("44÷2=22, 0", "48÷5=9, 3"),
    @("76÷8=9, 4", "79÷5=15, 4"),
    @("68÷4=17, 0", "90÷7=12, 6"),
    @("85÷4=21, 1", "25÷4=6, 1"),
    @("64÷8=8, 0", "86÷4=21, 2"),
    @("12÷5=2, 2", "75÷9=8, 3"),
    @("81÷2=40, 1", "34÷7=4, 6"),
    @("51÷9=5, 6", "66÷4=16, 2"),
    @("82÷8=10, 2", "19÷7=2, 5"),
    @("65÷8=8, 1", "12÷2=6, 0"),
    @("94÷8=11, 6", "18÷4=4, 2"),
    @("86÷3=28, 2", "22÷8=2, 6"),
    @("86÷2=43, 0", "46÷6=7, 4"),
    @("31÷6=5, 1", "42÷7=6, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
